function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws 'D2' '28.517.20'
Set-TextValue $ws 'E2' '  -1.74%  '
Set-TextValue $ws 'D3' '1.958.48'
Set-TextValue $ws 'E3' '  -0.23%  '
Set-TextValue $ws 'D4' '1.010'
Set-TextValue $ws 'E4' '  +0.24%  '
Set-TextValue $ws 'D5' '322.74'
Set-TextValue $ws 'E5' '  -1.25%  '
Set-TextValue $ws 'E6' '  +0.39%  '
Set-TextValue $ws 'D7' '0.4801'
Set-TextValue $ws 'E7' '  -3.75%  '
Set-TextValue $ws 'D8' '0.4090'
Set-TextValue $ws 'E8' '  -2.66%  '
Set-TextValue $ws 'D9' '54.04'
Set-TextValue $ws 'E9' '  +2.45%  '
Set-TextValue $ws 'D10' '0.08504'
Set-TextValue $ws 'E10' '  -6.96%  '
Set-TextValue $ws 'E11' '  -3.33%  '
Set-TextValue $ws 'D12' '22.50'
Set-TextValue $ws 'D13' '1.936.25'
Set-TextValue $ws 'E13' '  -1.02%  '
Set-TextValue $ws 'D14' '7.592'
Set-TextValue $ws 'E14' '  -3.19%  '
Set-TextValue $ws 'D15' '6.174'
Set-TextValue $ws 'D16' '1.013'
Set-TextValue $ws 'E16' '  +0.52%  '
Set-TextValue $ws 'D17' '90.72'
Set-TextValue $ws 'E17' '  -0.59%  '
Set-TextValue $ws 'D18' '0.00001070'
Set-TextValue $ws 'E18' '  -2.62%  '
Set-TextValue $ws 'D19' '0.06621'
Set-TextValue $ws 'E19' '  -1.06%  '
Set-TextValue $ws 'D20' '18.53'
Set-TextValue $ws 'E20' '  -3.38%  '
Set-TextValue $ws 'E21' '  +0.45%  '
Set-TextValue $ws 'D22' '5.847'
Set-TextValue $ws 'E22' '  -2.00%  '
Set-TextValue $ws 'D23' '28.512.80'
Set-TextValue $ws 'E23' '  -1.83%  '
Set-TextValue $ws 'E24' '  -4.58%  '
Set-TextValue $ws 'D25' '2.300'
Set-TextValue $ws 'E25' '  +0.74%  '
Set-TextValue $ws 'D26' '2.171.77'
Set-TextValue $ws 'E26' '  -1.23%  '
Set-TextValue $ws 'D27' '156.70'
Set-TextValue $ws 'E27' '  +0.37%  '
Set-TextValue $ws 'E28' '  -0.99%  '
Set-TextValue $ws 'D29' '2.179'
Set-TextValue $ws 'E29' '  -3.57%  '
Set-TextValue $ws 'D30' '5.841'
Set-TextValue $ws 'E30' '  -5.38%  '
Set-TextValue $ws 'D31' '124.43'
Set-TextValue $ws 'E31' '  -1.76%  '
Set-TextValue $ws 'D32' '0.9895'
Set-TextValue $ws 'E32' '  -4.54%  '
Set-TextValue $ws 'D33' '0.09682'
Set-TextValue $ws 'E33' '  -1.47%  '
Set-TextValue $ws 'D34' '1.455'
Set-TextValue $ws 'E34' '  -4.58%  '
Set-TextValue $ws 'B35' 'HuobiToken'
Set-TextValue $ws 'C35' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws 'D35' '3.693'
Set-TextValue $ws 'E35' '  +0.44%  '
Set-TextValue $ws 'B36' 'Filecoin'
Set-TextValue $ws 'C36' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws 'D36' '5.637'
Set-TextValue $ws 'E36' '  -2.13%  '
Set-TextValue $ws 'E37' '  +2.78%  '
Set-TextValue $ws 'D38' '0.02340'
Set-TextValue $ws 'E38' '  -3.07%  '
Set-TextValue $ws 'D39' '0.06204'
Set-TextValue $ws 'E39' '  -1.74%  '
Set-TextValue $ws 'D40' '1.255'
Set-TextValue $ws 'E40' '  -3.27%  '
Set-TextValue $ws 'D41' '0.6245'
Set-TextValue $ws 'E41' '  -2.88%  '
Set-TextValue $ws 'D42' '11.23'
Set-TextValue $ws 'E42' '  -1.65%  '
Set-TextValue $ws 'D43' '1.010'
Set-TextValue $ws 'E43' '  +0.42%  '
Set-TextValue $ws 'E44' '  -2.95%  '
Set-TextValue $ws 'D45' '1.344'
Set-TextValue $ws 'E45' '  +4.05%  '
Set-TextValue $ws 'D46' '0.5966'
Set-TextValue $ws 'E46' '  -4.02%  '
Set-TextValue $ws 'D47' '13.04'
Set-TextValue $ws 'E47' '  -1.80%  '
Set-TextValue $ws 'E48' '  -5.11%  '
Set-TextValue $ws 'D49' '3.411'
Set-TextValue $ws 'E49' '  -1.52%  '
Set-TextValue $ws 'D50' '0.06816'
Set-TextValue $ws 'E50' '  -1.81%  '
Set-TextValue $ws 'B51' 'BabyDogeCoin'
Set-TextValue $ws 'C51' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws 'D51' '0.00000000312'
Set-TextValue $ws 'E51' '  -6.16%  '
